$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 previously held an email address (shared string); the test data
# generator now leaves it blank while keeping the default cell style.
$ws.Range("D2").Value = ""

# D3 gets a freshly generated email address for this test run.
$ws.Range("D3").Value = "o67iEFMWXU@yopmail.com"

# O2 gets a new numeric-looking order/reference code. Typing it directly
# would make Excel coerce it to a Number (dropping the leading zeros) or
# stamp a "keep as text" number format onto the cell's style. Route the
# literal text through a formula result and Paste Special (values only)
# instead, which deposits it as a plain shared-string cell without
# touching styles.xml.
$ws.Range("O2").Formula = "=""000022472"""
$ws.Range("O2").Copy()
$ws.Range("O2").PasteSpecial(-4163)
$excel.CutCopyMode = $false
